# Updates cryptos list figures (price + 1h volume change) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "23.459.04"),
    @("E2", "  -0.26%  "),
    @("D3", "1.631.03"),
    @("E3", "  -0.51%  "),
    @("D4", "'1.001"),
    @("E4", "  +0.00%  "),
    @("E5", "  +0.06%  "),
    @("D6", "'305.13"),
    @("E6", "  -0.95%  "),
    @("D7", "'0.3759"),
    @("D8", "'0.3657"),
    @("E8", "  +0.04%  "),
    @("D9", "'51.86"),
    @("E9", "  -1.13%  "),
    @("D10", "'0.08225"),
    @("E10", "  +0.45%  "),
    @("D11", "'1.226"),
    @("E11", "  -3.61%  "),
    @("D12", "'1.002"),
    @("E12", "  +0.02%  "),
    @("D13", "'22.54"),
    @("E13", "  -2.00%  "),
    @("D14", "'6.567"),
    @("E14", "  -1.22%  "),
    @("D15", "'0.00001251"),
    @("E15", "  -2.21%  "),
    @("D16", "'7.261"),
    @("E16", "  -1.87%  "),
    @("D17", "1.631.71"),
    @("E17", "  -0.53%  "),
    @("D18", "'94.15"),
    @("E18", "  -0.65%  "),
    @("D19", "'0.06983"),
    @("E19", "  +0.56%  "),
    @("D20", "'17.78"),
    @("D21", "'6.462"),
    @("D22", "'1.002"),
    @("E22", "  +0.20%  "),
    @("D23", "'12.75"),
    @("E23", "  -0.60%  "),
    @("D24", "23.460.99"),
    @("E24", "  -0.28%  "),
    @("D25", "'3.166"),
    @("E25", "  +3.16%  "),
    @("D26", "'2.460"),
    @("E26", "  +1.53%  "),
    @("D27", "'21.42"),
    @("E27", "  +0.50%  "),
    @("D28", "'150.44"),
    @("E28", "  -0.64%  "),
    @("D29", "'5.315"),
    @("E29", "  -0.74%  "),
    @("D30", "'133.84"),
    @("E30", "  -1.38%  "),
    @("D31", "1.814.83"),
    @("E31", "  -0.37%  "),
    @("D32", "'2.257"),
    @("E32", "  -5.21%  "),
    @("D33", "'6.819"),
    @("E33", "  +0.12%  "),
    @("D34", "'1.018"),
    @("E34", "  +4.47%  "),
    @("D35", "'10.88"),
    @("E35", "  +5.26%  "),
    @("D36", "'0.02798"),
    @("E36", "  -1.03%  "),
    @("D37", "'0.2532"),
    @("E37", "  -0.84%  "),
    @("B38", "Stellar"),
    @("C38", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"),
    @("D38", "'0.08782"),
    @("E38", "  -1.30%  "),
    @("D39", "'0.07172"),
    @("E39", "  -2.74%  "),
    @("B40", "InternetComputer(DFINITY)"),
    @("C40", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"),
    @("D40", "'6.060"),
    @("E40", "  -2.09%  "),
    @("D41", "'0.7061"),
    @("E41", "  -0.75%  "),
    @("D42", "'1.350"),
    @("D43", "'16.38"),
    @("E43", "  +0.48%  "),
    @("D44", "'12.29"),
    @("D45", "'0.6570"),
    @("E45", "  +0.28%  "),
    @("D46", "'2.332"),
    @("E46", "  -0.60%  "),
    @("D47", "'1.000"),
    @("D48", "'3.991"),
    @("E48", "  -1.26%  "),
    @("D49", "'0.08023"),
    @("E49", "  +0.55%  "),
    @("D50", "'1.206"),
    @("E50", "  -0.57%  "),
    @("D51", "'125.58"),
    @("E51", "  -3.11%  ")
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
